# Insert a new bullet ("Instructor NPS score of 86.96, target of 50")
# immediately after the "Instructor engagement: 4.6/5" bullet, matching
# the formatting of other themed/size-20 bullets already present in the
# resume (e.g. the "Evangelized ..." bullet).

$d = $word.ActiveDocument

# Locate the "Instructor engagement" paragraph and remember its index so
# we can reliably address the paragraph that gets created right after it.
$target = $null
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Instructor engagement: 4.6/5*") {
        $target = $p
        $targetIndex = $i
        break
    }
    $i = $i + 1
}

if ($target -eq $null) {
    throw "Could not find the 'Instructor engagement: 4.6/5' paragraph"
}

# Create a new, empty paragraph right after it (inherits list/style
# context from Word the same way pressing Enter at the end of the line
# would), then stamp it with the exact OOXML the final paragraph should
# contain (flat-OPC wrapper, as required by Range.InsertXML).
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex + 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="30"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Instructor NPS score of 86.96, target of 50</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($xml)
